$wb = $excel.ActiveWorkbook

# --- Fix header typo "Percipitation_mm" -> "Precipitation_mm" on all sheets ---
foreach ($ws in $wb.Worksheets) {
    $ws.Range("F1").Value = "Precipitation_mm"
}

# --- sea_sunny: season spring -> summer for rows 2-49 ---
$wsSunny = $wb.Worksheets.Item("sea_sunny")
for ($r = 2; $r -le 49; $r++) {
    $wsSunny.Cells.Item($r, 5).Value = "summer"
}

# --- sea_windy & sea_cloudy: season autumn -> fall for rows 2-49 ---
$wsWindy = $wb.Worksheets.Item("sea_windy")
$wsCloudy = $wb.Worksheets.Item("sea_cloudy")
for ($r = 2; $r -le 49; $r++) {
    $wsWindy.Cells.Item($r, 5).Value = "fall"
    $wsCloudy.Cells.Item($r, 5).Value = "fall"
}

# --- sea_rainy: replace data rows 2-49 with new Denver weather data, season -> fall ---
$wsRainy = $wb.Worksheets.Item("sea_rainy")
$wsRainy.Range("A2").Value = 44498
$wsRainy.Range("B2").Value = 9.699999999999999
$wsRainy.Range("C2").Value = 0
$wsRainy.Range("D2").Value = 0.4
$wsRainy.Range("E2").Value = "fall"
$wsRainy.Range("F2").Value = 1.905
$wsRainy.Range("A3").Value = 44498.04166666666
$wsRainy.Range("B3").Value = 9.300000000000001
$wsRainy.Range("C3").Value = 0
$wsRainy.Range("D3").Value = 0.4
$wsRainy.Range("E3").Value = "fall"
$wsRainy.Range("F3").Value = 0.801076923076923
$wsRainy.Range("A4").Value = 44498.08333333334
$wsRainy.Range("B4").Value = 9
$wsRainy.Range("C4").Value = 0
$wsRainy.Range("D4").Value = 0.5
$wsRainy.Range("E4").Value = "fall"
$wsRainy.Range("F4").Value = 1.733176470588235
$wsRainy.Range("A5").Value = 44498.125
$wsRainy.Range("B5").Value = 8.800000000000001
$wsRainy.Range("C5").Value = 0
$wsRainy.Range("D5").Value = 0.5
$wsRainy.Range("E5").Value = "fall"
$wsRainy.Range("F5").Value = 1.23825
$wsRainy.Range("A6").Value = 44498.16666666666
$wsRainy.Range("B6").Value = 8.5
$wsRainy.Range("C6").Value = 0
$wsRainy.Range("D6").Value = 0.4
$wsRainy.Range("E6").Value = "fall"
$wsRainy.Range("F6").Value = 1.172307692307692
$wsRainy.Range("A7").Value = 44498.20833333334
$wsRainy.Range("B7").Value = 8.4
$wsRainy.Range("C7").Value = 0
$wsRainy.Range("D7").Value = 0.5
$wsRainy.Range("E7").Value = "fall"
$wsRainy.Range("F7").Value = 0.7281333333333334
$wsRainy.Range("A8").Value = 44498.25
$wsRainy.Range("B8").Value = 8.4
$wsRainy.Range("C8").Value = 0
$wsRainy.Range("D8").Value = 0.5
$wsRainy.Range("E8").Value = "fall"
$wsRainy.Range("F8").Value = 0.7619999999999999
$wsRainy.Range("A9").Value = 44498.29166666666
$wsRainy.Range("B9").Value = 8.1
$wsRainy.Range("C9").Value = 0
$wsRainy.Range("D9").Value = 0.6000000000000001
$wsRainy.Range("E9").Value = "fall"
$wsRainy.Range("F9").Value = 0.4535714285714286
$wsRainy.Range("A10").Value = 44498.33333333334
$wsRainy.Range("B10").Value = 7.800000000000001
$wsRainy.Range("C10").Value = 2
$wsRainy.Range("D10").Value = 0.7000000000000001
$wsRainy.Range("E10").Value = "fall"
$wsRainy.Range("F10").Value = 0.3321538461538461
$wsRainy.Range("A11").Value = 44498.375
$wsRainy.Range("B11").Value = 9
$wsRainy.Range("C11").Value = 18
$wsRainy.Range("D11").Value = 0.4
$wsRainy.Range("E11").Value = "fall"
$wsRainy.Range("F11").Value = 0.1953846153846154
$wsRainy.Range("A12").Value = 44498.41666666666
$wsRainy.Range("B12").Value = 10.5
$wsRainy.Range("C12").Value = 39
$wsRainy.Range("D12").Value = 0.5
$wsRainy.Range("E12").Value = "fall"
$wsRainy.Range("F12").Value = 0.1693333333333333
$wsRainy.Range("A13").Value = 44498.45833333334
$wsRainy.Range("B13").Value = 11.2
$wsRainy.Range("C13").Value = 64
$wsRainy.Range("D13").Value = 0.6000000000000001
$wsRainy.Range("E13").Value = "fall"
$wsRainy.Range("F13").Value = 0
$wsRainy.Range("A14").Value = 44498.5
$wsRainy.Range("B14").Value = 12
$wsRainy.Range("C14").Value = 122
$wsRainy.Range("D14").Value = 0.6000000000000001
$wsRainy.Range("E14").Value = "fall"
$wsRainy.Range("F14").Value = 0
$wsRainy.Range("A15").Value = 44498.54166666666
$wsRainy.Range("B15").Value = 12.2
$wsRainy.Range("C15").Value = 50
$wsRainy.Range("D15").Value = 0.6000000000000001
$wsRainy.Range("E15").Value = "fall"
$wsRainy.Range("F15").Value = 0.001923076923076923
$wsRainy.Range("A16").Value = 44498.58333333334
$wsRainy.Range("B16").Value = 12
$wsRainy.Range("C16").Value = 45
$wsRainy.Range("D16").Value = 0.6000000000000001
$wsRainy.Range("E16").Value = "fall"
$wsRainy.Range("F16").Value = 0.003571428571428572
$wsRainy.Range("A17").Value = 44498.625
$wsRainy.Range("B17").Value = 11.8
$wsRainy.Range("C17").Value = 34
$wsRainy.Range("D17").Value = 0.6000000000000001
$wsRainy.Range("E17").Value = "fall"
$wsRainy.Range("F17").Value = 0
$wsRainy.Range("A18").Value = 44498.66666666666
$wsRainy.Range("B18").Value = 11
$wsRainy.Range("C18").Value = 28
$wsRainy.Range("D18").Value = 0.6000000000000001
$wsRainy.Range("E18").Value = "fall"
$wsRainy.Range("F18").Value = 0
$wsRainy.Range("A19").Value = 44498.70833333334
$wsRainy.Range("B19").Value = 9.199999999999999
$wsRainy.Range("C19").Value = 106
$wsRainy.Range("D19").Value = 0.6000000000000001
$wsRainy.Range("E19").Value = "fall"
$wsRainy.Range("F19").Value = 0
$wsRainy.Range("A20").Value = 44498.75
$wsRainy.Range("B20").Value = 6.7
$wsRainy.Range("C20").Value = 0
$wsRainy.Range("D20").Value = 0.8
$wsRainy.Range("E20").Value = "fall"
$wsRainy.Range("F20").Value = 0
$wsRainy.Range("A21").Value = 44498.79166666666
$wsRainy.Range("B21").Value = 6.100000000000001
$wsRainy.Range("C21").Value = 0
$wsRainy.Range("D21").Value = 0.9
$wsRainy.Range("E21").Value = "fall"
$wsRainy.Range("F21").Value = 0
$wsRainy.Range("A22").Value = 44498.83333333334
$wsRainy.Range("B22").Value = 5.5
$wsRainy.Range("C22").Value = 0
$wsRainy.Range("D22").Value = 0.9
$wsRainy.Range("E22").Value = "fall"
$wsRainy.Range("F22").Value = 0
$wsRainy.Range("A23").Value = 44498.875
$wsRainy.Range("B23").Value = 4.9
$wsRainy.Range("C23").Value = 0
$wsRainy.Range("D23").Value = 0.9
$wsRainy.Range("E23").Value = "fall"
$wsRainy.Range("F23").Value = 0
$wsRainy.Range("A24").Value = 44498.91666666666
$wsRainy.Range("B24").Value = 4.600000000000001
$wsRainy.Range("C24").Value = 0
$wsRainy.Range("D24").Value = 1
$wsRainy.Range("E24").Value = "fall"
$wsRainy.Range("F24").Value = 0
$wsRainy.Range("A25").Value = 44498.95833333334
$wsRainy.Range("B25").Value = 4.3
$wsRainy.Range("C25").Value = 0
$wsRainy.Range("D25").Value = 1
$wsRainy.Range("E25").Value = "fall"
$wsRainy.Range("F25").Value = 0
$wsRainy.Range("A26").Value = 44499
$wsRainy.Range("B26").Value = 3.9
$wsRainy.Range("C26").Value = 0
$wsRainy.Range("D26").Value = 0.9
$wsRainy.Range("E26").Value = "fall"
$wsRainy.Range("F26").Value = 0
$wsRainy.Range("A27").Value = 44499.04166666666
$wsRainy.Range("B27").Value = 3.7
$wsRainy.Range("C27").Value = 0
$wsRainy.Range("D27").Value = 0.8
$wsRainy.Range("E27").Value = "fall"
$wsRainy.Range("F27").Value = 0
$wsRainy.Range("A28").Value = 44499.08333333334
$wsRainy.Range("B28").Value = 3.9
$wsRainy.Range("C28").Value = 0
$wsRainy.Range("D28").Value = 0.7000000000000001
$wsRainy.Range("E28").Value = "fall"
$wsRainy.Range("F28").Value = 0
$wsRainy.Range("A29").Value = 44499.125
$wsRainy.Range("B29").Value = 3.8
$wsRainy.Range("C29").Value = 0
$wsRainy.Range("D29").Value = 0.7000000000000001
$wsRainy.Range("E29").Value = "fall"
$wsRainy.Range("F29").Value = 0
$wsRainy.Range("A30").Value = 44499.16666666666
$wsRainy.Range("B30").Value = 3.7
$wsRainy.Range("C30").Value = 0
$wsRainy.Range("D30").Value = 0.8
$wsRainy.Range("E30").Value = "fall"
$wsRainy.Range("F30").Value = 0
$wsRainy.Range("A31").Value = 44499.20833333334
$wsRainy.Range("B31").Value = 2.9
$wsRainy.Range("C31").Value = 0
$wsRainy.Range("D31").Value = 0.7000000000000001
$wsRainy.Range("E31").Value = "fall"
$wsRainy.Range("F31").Value = 0
$wsRainy.Range("A32").Value = 44499.25
$wsRainy.Range("B32").Value = 2.7
$wsRainy.Range("C32").Value = 0
$wsRainy.Range("D32").Value = 0.7000000000000001
$wsRainy.Range("E32").Value = "fall"
$wsRainy.Range("F32").Value = 0
$wsRainy.Range("A33").Value = 44499.29166666666
$wsRainy.Range("B33").Value = 2.2
$wsRainy.Range("C33").Value = 0
$wsRainy.Range("D33").Value = 0.6000000000000001
$wsRainy.Range("E33").Value = "fall"
$wsRainy.Range("F33").Value = 0
$wsRainy.Range("A34").Value = 44499.33333333334
$wsRainy.Range("B34").Value = 2.4
$wsRainy.Range("C34").Value = 14
$wsRainy.Range("D34").Value = 0.6000000000000001
$wsRainy.Range("E34").Value = "fall"
$wsRainy.Range("F34").Value = 0
$wsRainy.Range("A35").Value = 44499.375
$wsRainy.Range("B35").Value = 4.5
$wsRainy.Range("C35").Value = 143
$wsRainy.Range("D35").Value = 0.7000000000000001
$wsRainy.Range("E35").Value = "fall"
$wsRainy.Range("F35").Value = 0
$wsRainy.Range("A36").Value = 44499.41666666666
$wsRainy.Range("B36").Value = 6.600000000000001
$wsRainy.Range("C36").Value = 289
$wsRainy.Range("D36").Value = 0.8
$wsRainy.Range("E36").Value = "fall"
$wsRainy.Range("F36").Value = 0
$wsRainy.Range("A37").Value = 44499.45833333334
$wsRainy.Range("B37").Value = 8.800000000000001
$wsRainy.Range("C37").Value = 407
$wsRainy.Range("D37").Value = 0.9
$wsRainy.Range("E37").Value = "fall"
$wsRainy.Range("F37").Value = 0
$wsRainy.Range("A38").Value = 44499.5
$wsRainy.Range("B38").Value = 10.6
$wsRainy.Range("C38").Value = 476
$wsRainy.Range("D38").Value = 0.9
$wsRainy.Range("E38").Value = "fall"
$wsRainy.Range("F38").Value = 0
$wsRainy.Range("A39").Value = 44499.54166666666
$wsRainy.Range("B39").Value = 11.7
$wsRainy.Range("C39").Value = 496
$wsRainy.Range("D39").Value = 1
$wsRainy.Range("E39").Value = "fall"
$wsRainy.Range("F39").Value = 0
$wsRainy.Range("A40").Value = 44499.58333333334
$wsRainy.Range("B40").Value = 12.1
$wsRainy.Range("C40").Value = 463
$wsRainy.Range("D40").Value = 1
$wsRainy.Range("E40").Value = "fall"
$wsRainy.Range("F40").Value = 0
$wsRainy.Range("A41").Value = 44499.625
$wsRainy.Range("B41").Value = 12
$wsRainy.Range("C41").Value = 379
$wsRainy.Range("D41").Value = 1
$wsRainy.Range("E41").Value = "fall"
$wsRainy.Range("F41").Value = 0
$wsRainy.Range("A42").Value = 44499.66666666666
$wsRainy.Range("B42").Value = 11.2
$wsRainy.Range("C42").Value = 253
$wsRainy.Range("D42").Value = 0.9
$wsRainy.Range("E42").Value = "fall"
$wsRainy.Range("F42").Value = 0
$wsRainy.Range("A43").Value = 44499.70833333334
$wsRainy.Range("B43").Value = 8.800000000000001
$wsRainy.Range("C43").Value = 106
$wsRainy.Range("D43").Value = 0.8
$wsRainy.Range("E43").Value = "fall"
$wsRainy.Range("F43").Value = 0
$wsRainy.Range("A44").Value = 44499.75
$wsRainy.Range("B44").Value = 6
$wsRainy.Range("C44").Value = 0
$wsRainy.Range("D44").Value = 0.9
$wsRainy.Range("E44").Value = "fall"
$wsRainy.Range("F44").Value = 0
$wsRainy.Range("A45").Value = 44499.79166666666
$wsRainy.Range("B45").Value = 5.5
$wsRainy.Range("C45").Value = 0
$wsRainy.Range("D45").Value = 0.9
$wsRainy.Range("E45").Value = "fall"
$wsRainy.Range("F45").Value = 0
$wsRainy.Range("A46").Value = 44499.83333333334
$wsRainy.Range("B46").Value = 5
$wsRainy.Range("C46").Value = 0
$wsRainy.Range("D46").Value = 0.8
$wsRainy.Range("E46").Value = "fall"
$wsRainy.Range("F46").Value = 0
$wsRainy.Range("A47").Value = 44499.875
$wsRainy.Range("B47").Value = 4.5
$wsRainy.Range("C47").Value = 0
$wsRainy.Range("D47").Value = 0.8
$wsRainy.Range("E47").Value = "fall"
$wsRainy.Range("F47").Value = 0
$wsRainy.Range("A48").Value = 44499.91666666666
$wsRainy.Range("B48").Value = 4.4
$wsRainy.Range("C48").Value = 0
$wsRainy.Range("D48").Value = 0.8
$wsRainy.Range("E48").Value = "fall"
$wsRainy.Range("F48").Value = 0
$wsRainy.Range("A49").Value = 44499.95833333334
$wsRainy.Range("B49").Value = 4.5
$wsRainy.Range("C49").Value = 0
$wsRainy.Range("D49").Value = 1
$wsRainy.Range("E49").Value = "fall"
$wsRainy.Range("F49").Value = 0
